# Regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# Updates column G ("K") values for rows 2-38 on the active sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 1
    3  = 1
    4  = 7
    5  = 9
    6  = 3
    7  = 7
    8  = 2
    9  = 2
    10 = 9
    11 = 3
    12 = 11
    13 = 3
    14 = 4
    15 = 4
    16 = 0
    17 = 5
    18 = 5
    19 = 2
    20 = 2
    21 = 4
    22 = 6
    23 = 4
    24 = 3
    25 = 4
    26 = 7
    27 = 7
    28 = 1
    29 = 3
    30 = 3
    31 = 4
    32 = 3
    33 = 6
    34 = 10
    35 = 4
    36 = 4
    37 = 3
    38 = 3
}

foreach ($row in $newK.Keys) {
    $ws.Range("G$row").Value = $newK[$row]
}
